$d = $word.ActiveDocument

$replacements = @(
    @{old="829÷6=138, 1"; new="532÷9=59, 1"},
    @{old="603÷6=100, 3"; new="909÷7=129, 6"},
    @{old="336÷3=112, 0"; new="471÷7=67, 2"},
    @{old="591÷2=295, 1"; new="880÷8=110, 0"},
    @{old="849÷4=212, 1"; new="835÷7=119, 2"},
    @{old="141÷9=15, 6"; new="327÷2=163, 1"},
    @{old="651÷5=130, 1"; new="649÷2=324, 1"},
    @{old="506÷8=63, 2"; new="948÷9=105, 3"},
    @{old="508÷6=84, 4"; new="214÷2=107, 0"},
    @{old="677÷4=169, 1"; new="768÷7=109, 5"},
    @{old="489÷5=97, 4"; new="279÷3=93, 0"},
    @{old="213÷8=26, 5"; new="712÷8=89, 0"},
    @{old="707÷2=353, 1"; new="476÷7=68, 0"},
    @{old="679÷6=113, 1"; new="646÷9=71, 7"},
    @{old="601÷8=75, 1"; new="530÷3=176, 2"},
    @{old="819÷5=163, 4"; new="214÷9=23, 7"},
    @{old="146÷9=16, 2"; new="148÷9=16, 4"},
    @{old="533÷2=266, 1"; new="750÷6=125, 0"},
    @{old="851÷9=94, 5"; new="248÷3=82, 2"},
    @{old="547÷2=273, 1"; new="978÷4=244, 2"},
    @{old="377÷4=94, 1"; new="914÷5=182, 4"},
    @{old="978÷3=326, 0"; new="884÷2=442, 0"},
    @{old="814÷9=90, 4"; new="193÷2=96, 1"},
    @{old="402÷5=80, 2"; new="276÷3=92, 0"},
    @{old="962÷2=481, 0"; new="288÷5=57, 3"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
